$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row additions (F1:I1), matching existing bold header style ---
$ws.Range("F1").Value = "fruit_time"
$ws.Range("G1").Value = "locality"
$ws.Range("H1").Value = "flower_time"
$ws.Range("I1").Value = "flower_locality"
$ws.Range("F1:I1").Font.Bold = $true

# --- Column E width -> 16 ---
$ws.Columns.Item(5).ColumnWidth = 15.15

# --- Data rows: new flowering-time columns (F-I) + a few corrected C values ---
$ws.Range("F2").Value = "early_summer"
$ws.Range("G2").Value = "southeast"
$ws.Range("H2").Value = 4
$ws.Range("F3").Value = 11.5
$ws.Range("G3").Value = "oregon?"
$ws.Range("H3").Value = 4.5
$ws.Range("F4").Value = 9.5
$ws.Range("G4").Value = "New england"
$ws.Range("H4").Value = 4.5
$ws.Range("I4").Value = "range ave"
$ws.Range("C5").Value = "wind"
$ws.Range("F5").Value = 8.5
$ws.Range("H5").Value = 4
$ws.Range("C6").Value = "wind"
$ws.Range("F6").Value = 5.5
$ws.Range("H6").Value = 4
$ws.Range("C7").Value = "wind"
$ws.Range("F7").Value = 5
$ws.Range("H7").Value = 3.5
$ws.Range("C8").Value = "insect"
$ws.Range("F8").Value = 10.5
$ws.Range("H8").Value = 5.5
$ws.Range("F9").Value = 9
$ws.Range("H9").Value = 4
$ws.Range("C10").Value = "insect"
$ws.Range("F10").Value = 9.5
$ws.Range("H10").Value = 4
$ws.Range("C11").Value = "insect"
$ws.Range("F11").Value = 9
$ws.Range("H11").Value = 5
$ws.Range("C12").Value = "wind"
$ws.Range("F12").Value = 8.5
$ws.Range("H12").Value = 5.5
$ws.Range("F13").Value = 10
$ws.Range("H13").Value = 4.5
$ws.Range("F14").Value = "late_spring/early_summer"
$ws.Range("H14").Value = 4.5
$ws.Range("F15").Value = "NA"
$ws.Range("H15").Value = 4.5
$ws.Range("I15").Value = "north"
$ws.Range("F16").Value = 11
$ws.Range("H16").Value = 4.5
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = 5
$ws.Range("F18").Value = 10.5
$ws.Range("H18").Value = 4.5
$ws.Range("F19").Value = 10.5
$ws.Range("H19").Value = 6
$ws.Range("I19").Value = "north"
$ws.Range("F20").Value = "NA"
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = "north"
$ws.Range("F21").Value = "as late as midwinter"
$ws.Range("H21").Value = 4
$ws.Range("F22").Value = "fall/winter"
$ws.Range("H22").Value = 4
$ws.Range("F23").Value = "NA"
$ws.Range("H23").Value = 5
$ws.Range("I23").Value = "north"
$ws.Range("F24").Value = 10
$ws.Range("H24").Value = 4.5
$ws.Range("F25").Value = 10
$ws.Range("H25").Value = 4.5
$ws.Range("F26").Value = 10
$ws.Range("H26").Value = 4.5
$ws.Range("F27").Value = 8.5
$ws.Range("H27").Value = 4.5
$ws.Range("F28").Value = 8.5
$ws.Range("H28").Value = 5.5
$ws.Range("F29").Value = 9.5
$ws.Range("H29").Value = 4.5
$ws.Range("I29").Value = "north"
$ws.Range("F30").Value = 11
$ws.Range("H30").Value = 4.5
$ws.Range("F31").Value = 10.5
$ws.Range("H31").Value = 6
$ws.Range("I31").Value = "north"
$ws.Range("F32").Value = 10.5
$ws.Range("H32").Value = 6.5
$ws.Range("F33").Value = "winter"
$ws.Range("H33").Value = 5
$ws.Range("F34").Value = 10
$ws.Range("H34").Value = 6
$ws.Range("I34").Value = "north"
$ws.Range("F35").Value = 9
$ws.Range("H35").Value = 5.5
$ws.Range("F36").Value = 9.5
$ws.Range("H36").Value = 5
$ws.Range("F37").Value = 10.5
$ws.Range("H37").Value = 3.5
$ws.Range("F38").Value = 11.5
$ws.Range("H38").Value = 4
$ws.Range("F39").Value = 9
$ws.Range("G39").Value = "range"
$ws.Range("H39").Value = 6
$ws.Range("I39").Value = "north"
$ws.Range("F40").Value = 6.5
$ws.Range("G40").Value = "north"
$ws.Range("H40").Value = 3
$ws.Range("F41").Value = 5.5
$ws.Range("H41").Value = 4.5
$ws.Range("F42").Value = 5.5
$ws.Range("H42").Value = 4
$ws.Range("F43").Value = 5.5
$ws.Range("G43").Value = "range"
$ws.Range("H43").Value = 4.5
$ws.Range("I43").Value = "range"
$ws.Range("F44").Value = 5.5
$ws.Range("H44").Value = 4
$ws.Range("C45").Value = "insect"
$ws.Range("F45").Value = 7
$ws.Range("H45").Value = 5
$ws.Range("F46").Value = 8.5
$ws.Range("G46").Value = "allegheny"
$ws.Range("H46").Value = 5
$ws.Range("F47").Value = 9.5
$ws.Range("H47").Value = 4
$ws.Range("F48").Value = 9.5
$ws.Range("H48").Value = 5.5
$ws.Range("F49").Value = 9.5
$ws.Range("H49").Value = 6.5
$ws.Range("I49").Value = "north"
$ws.Range("F50").Value = 21
$ws.Range("H50").Value = 4.5
$ws.Range("F51").Value = 21
$ws.Range("H51").Value = 4.5
$ws.Range("F52").Value = 22
$ws.Range("H52").Value = "NA"
$ws.Range("F53").Value = 9.5
$ws.Range("H53").Value = 4.5
$ws.Range("F54").Value = 21
$ws.Range("H54").Value = 2.5
$ws.Range("F55").Value = 9.5
$ws.Range("H55").Value = 4.5
$ws.Range("F56").Value = 9.5
$ws.Range("H56").Value = 4.5
$ws.Range("F57").Value = 21
$ws.Range("H57").Value = "NA"
$ws.Range("F58").Value = 11.5
$ws.Range("H58").Value = 3.5
$ws.Range("F59").Value = 21
$ws.Range("H59").Value = 3.5
$ws.Range("F60").Value = 10
$ws.Range("H60").Value = 5
$ws.Range("I60").Value = "north"
$ws.Range("F61").Value = 9.5
$ws.Range("H61").Value = "NA"
$ws.Range("F62").Value = 12.5
$ws.Range("H62").Value = 5.5
$ws.Range("F63").Value = 6.5
$ws.Range("H63").Value = 4
$ws.Range("I63").Value = "range"
$ws.Range("F64").Value = 8.5
$ws.Range("H64").Value = 3.5
$ws.Range("F65").Value = 9.5
$ws.Range("H65").Value = 6
$ws.Range("F66").Value = 4
$ws.Range("H66").Value = 3.5
$ws.Range("F67").Value = 4.5
$ws.Range("G67").Value = "range"
$ws.Range("H67").Value = 3.5
$ws.Range("I67").Value = "range"
$ws.Range("F68").Value = 5.5
$ws.Range("H68").Value = 4
$ws.Range("F69").Value = 5
$ws.Range("H69").Value = 3.5

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Final selection to match author view ---
$ws.Range("H70").Select()

Write-Host "edit complete"
